# Applies the changes described in the commit diff:
#  - Bumps the "Date" metadata value
#  - Appends "|4.0.1" FHIR version markers to ValueSet URLs and Reference(...)/
#    Quantity {...} type strings on the Elements sheet
#  - Updates the (auto) column widths for columns K and Z on the Elements sheet
#    to reflect the widened cell content

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet: bump the generation Date
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-09-23T14:10:57+00:00"

# ---------------------------------------------------------------------------
# Elements sheet: append |4.0.1 to ValueSet URLs and Reference(...)/Quantity types
# ---------------------------------------------------------------------------
$el = $wb.Worksheets.Item("Elements")

# Binding Value Set (column Z) URLs
# (Some of these values are repeated lower in the sheet - rows 47/49/50 reuse
#  the same text as rows 16/24/25, so every occurrence must be updated
#  individually since the COM layer does not re-dedupe shared strings.)
$el.Range("Z6").Value  = "http://hl7.org/fhir/ValueSet/languages|4.0.1"
$el.Range("Z15").Value = "http://hl7.org/fhir/ValueSet/observation-category|4.0.1"
$el.Range("Z16").Value = "http://hl7.org/fhir/ValueSet/observation-codes|4.0.1"
$el.Range("Z24").Value = "http://hl7.org/fhir/ValueSet/data-absent-reason|4.0.1"
$el.Range("Z25").Value = "http://hl7.org/fhir/ValueSet/observation-interpretation|4.0.1"
$el.Range("Z27").Value = "http://hl7.org/fhir/ValueSet/body-site|4.0.1"
$el.Range("Z28").Value = "http://hl7.org/fhir/ValueSet/observation-methods|4.0.1"
$el.Range("Z37").Value = "http://hl7.org/fhir/ValueSet/referencerange-meaning|4.0.1"
$el.Range("Z38").Value = "http://hl7.org/fhir/ValueSet/referencerange-appliesto|4.0.1"
$el.Range("Z47").Value = "http://hl7.org/fhir/ValueSet/observation-codes|4.0.1"
$el.Range("Z49").Value = "http://hl7.org/fhir/ValueSet/data-absent-reason|4.0.1"
$el.Range("Z50").Value = "http://hl7.org/fhir/ValueSet/observation-interpretation|4.0.1"

# Type(s) (column K) Reference(...) / Quantity {...} values
$el.Range("K12").Value = "Reference(CarePlan|4.0.1|DeviceRequest|4.0.1|ImmunizationRecommendation|4.0.1|MedicationRequest|4.0.1|NutritionOrder|4.0.1|ServiceRequest|4.0.1)`n"
$el.Range("K13").Value = "Reference(MedicationAdministration|4.0.1|MedicationDispense|4.0.1|MedicationStatement|4.0.1|Procedure|4.0.1|Immunization|4.0.1|ImagingStudy|4.0.1)`n"
$el.Range("K18").Value = "Reference(Resource|4.0.1)`n"
$el.Range("K19").Value = "Reference(Encounter|4.0.1)`n"
$el.Range("K22").Value = "Reference(Practitioner|4.0.1|PractitionerRole|4.0.1|Organization|4.0.1|CareTeam|4.0.1|Patient|4.0.1|RelatedPerson|4.0.1)`n"
$el.Range("K29").Value = "Reference(Specimen|4.0.1)`n"
$el.Range("K30").Value = "Reference(Device|4.0.1|DeviceMetric|4.0.1)`n"
$el.Range("K35").Value = "Quantity {SimpleQuantity|4.0.1}`n"
$el.Range("K36").Value = "Quantity {SimpleQuantity|4.0.1}`n"
$el.Range("K41").Value = "Reference(Observation|4.0.1|QuestionnaireResponse|4.0.1|MolecularSequence|4.0.1)`n"
$el.Range("K42").Value = "Reference(DocumentReference|4.0.1|ImagingStudy|4.0.1|Media|4.0.1|QuestionnaireResponse|4.0.1|Observation|4.0.1|MolecularSequence|4.0.1)`n"

# ---------------------------------------------------------------------------
# Resize columns K (Type(s)) and Z (Binding Value Set) to reflect the new,
# wider, content (bestFit widths grew from 95.26171875 -> 122.43359375 and
# 43.21484375 -> 47.7421875 respectively).
# ---------------------------------------------------------------------------
$el.Columns.Item(11).ColumnWidth = 121.66666666666667
$el.Columns.Item(26).ColumnWidth = 46.833333333333336
